# Reorder the last 7 rows of the comorbidity table (table 1 in the document).
#
# Before (Word 1-based row index within Table 1):
#   24 Chronic pulmonary disease    | 1.76; < 0.01 | 1.66; < 0.01 | 0.98; 1.000
#   25 Metastatic solid tumour      | 1.39; 0.080  | 1.89; < 0.01 | 1.92; < 0.01
#   26 Rheumatic disease            | 1.56; 0.336  | 0.99; 1.000  | 0.79; 1.000
#   27 Peripheral vascular disease  | 2.49; < 0.01 | 1.54; < 0.01 | 1.04; 0.831
#   28 Renal disease                | 5.52; < 0.01 | 2.32; < 0.01 | 2.06; < 0.01
#   29 Mild liver disease           | 4.53; < 0.01 | 1.03; 0.833  | 1.11; 0.663
#   30 Congestive heart failure     | 5.53; < 0.01 | 3.73; < 0.01 | 1.67; 0.063
#
# After (rows 24-27 below replace rows 24-26; the four disease rows that used
# to sit before "Renal disease" now sit after "Congestive heart failure", in
# reverse order):
#   24 Renal disease                | 5.52; < 0.01 | 2.32; < 0.01 | 2.06; < 0.01
#   25 Mild liver disease           | 4.53; < 0.01 | 1.03; 0.833  | 1.11; 0.663
#   26 Congestive heart failure     | 5.53; < 0.01 | 3.73; < 0.01 | 1.67; 0.063
#   27 Peripheral vascular disease  | 2.49; < 0.01 | 1.54; < 0.01 | 1.04; 0.831
#   28 Metastatic solid tumour      | 1.39; 0.080  | 1.89; < 0.01 | 1.92; < 0.01
#   29 Rheumatic disease            | 1.56; 0.336  | 0.99; 1.000  | 0.79; 1.000
#   30 Chronic pulmonary disease    | 1.76; < 0.01 | 1.66; < 0.01 | 0.98; 1.000
#
# Rather than physically deleting/inserting rows (which would stamp fresh
# rsid/paraId attributes on brand new <w:tr> elements), we keep the existing
# seven <w:tr> elements in place and simply rewrite their cell text so the
# table ends up in the target order -- row 23 ("Any malignancy...") and
# everything above it is untouched.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$newRowsData = @(
    @("Renal disease", "5.52; < 0.01", "2.32; < 0.01", "2.06; < 0.01"),
    @("Mild liver disease", "4.53; < 0.01", "1.03; 0.833", "1.11; 0.663"),
    @("Congestive heart failure", "5.53; < 0.01", "3.73; < 0.01", "1.67; 0.063"),
    @("Peripheral vascular disease", "2.49; < 0.01", "1.54; < 0.01", "1.04; 0.831"),
    @("Metastatic solid tumour", "1.39; 0.080", "1.89; < 0.01", "1.92; < 0.01"),
    @("Rheumatic disease", "1.56; 0.336", "0.99; 1.000", "0.79; 1.000"),
    @("Chronic pulmonary disease", "1.76; < 0.01", "1.66; < 0.01", "0.98; 1.000")
)

$startRow = 24
for ($i = 0; $i -lt $newRowsData.Count; $i++) {
    $row = $t.Rows.Item($startRow + $i)
    $data = $newRowsData[$i]
    for ($c = 1; $c -le 4; $c++) {
        $row.Cells.Item($c).Range.Text = $data[$c - 1]
    }
}
